$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row with two new columns (P, Q), copying the header style from O1 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$headerArr = New-Object "object[,]" 1,2
$headerArr[0,0] = 14
$headerArr[0,1] = 15
$ws.Range("P1:Q1").Value = $headerArr

# --- Column H no longer has any data for rows 2:25 ---
$ws.Range("H2:H25").ClearContents()

# --- Updated simulation results for rows 2:25 ---
$BArr = New-Object "object[,]" 24,1
$BArr[0,0] = 24.74038215810859
$BArr[1,0] = 23.11024366375046
$BArr[2,0] = 22.05031507671086
$BArr[3,0] = 21.60325755213291
$BArr[4,0] = 21.5278557199714
$BArr[5,0] = 22.0436383199598
$BArr[6,0] = 24.18965281943812
$BArr[7,0] = 27.93442201274501
$BArr[8,0] = 30.40772020325911
$BArr[9,0] = 31.47310337544135
$BArr[10,0] = 31.86865553312921
$BArr[11,0] = 31.78392708364931
$BArr[12,0] = 31.5058513439434
$BArr[13,0] = 31.33425656028081
$BArr[14,0] = 30.33535654883388
$BArr[15,0] = 29.70756200589915
$BArr[16,0] = 29.34140133159286
$BArr[17,0] = 29.21609886079073
$BArr[18,0] = 29.77500266814817
$BArr[19,0] = 31.58724288809329
$BArr[20,0] = 32.72367984419313
$BArr[21,0] = 32.12216423825537
$BArr[22,0] = 29.74554229594099
$BArr[23,0] = 26.96950157942045
$ws.Range("B2:B25").Value = $BArr

$CArr = New-Object "object[,]" 24,1
$CArr[0,0] = 19.5636648860633
$CArr[1,0] = 18.41447534527523
$CArr[2,0] = 17.67850506758756
$CArr[3,0] = 17.38639751806545
$CArr[4,0] = 17.35625511068535
$CArr[5,0] = 17.72496012806476
$CArr[6,0] = 19.23605219913242
$CArr[7,0] = 21.87521780433195
$CArr[8,0] = 23.68702372883381
$CArr[9,0] = 24.54120525844531
$CArr[10,0] = 24.82606913706412
$CArr[11,0] = 24.75824459148323
$CArr[12,0] = 24.56187542589756
$CArr[13,0] = 24.45434137348789
$CArr[14,0] = 23.74521113548175
$CArr[15,0] = 23.3021713885696
$CArr[16,0] = 23.00866558759142
$CArr[17,0] = 22.9362724434703
$CArr[18,0] = 23.34686654937683
$CArr[19,0] = 24.65724537206995
$CArr[20,0] = 25.439988057616
$CArr[21,0] = 24.97970961732044
$CArr[22,0] = 23.25189885899947
$CArr[23,0] = 21.26521956912221
$ws.Range("C2:C25").Value = $CArr

$DArr = New-Object "object[,]" 24,1
$DArr[0,0] = 4.030612846141484
$DArr[1,0] = 3.90889586098775
$DArr[2,0] = 3.832792527754622
$DArr[3,0] = 3.8038479697238
$DArr[4,0] = 3.801949493575898
$DArr[5,0] = 3.840354864596208
$DArr[6,0] = 3.999276235006756
$DArr[7,0] = 4.28519168171126
$DArr[8,0] = 4.490633029677658
$DArr[9,0] = 4.593864647135679
$DArr[10,0] = 4.626297573667372
$DArr[11,0] = 4.618061515611747
$DArr[12,0] = 4.595999927109
$DArr[13,0] = 4.584954841394548
$DArr[14,0] = 4.50515019928633
$DArr[15,0] = 4.455623753089434
$DArr[16,0] = 4.420389858226298
$DArr[17,0] = 4.413465966819228
$DArr[18,0] = 4.460406054748779
$DArr[19,0] = 4.609596873267936
$DArr[20,0] = 4.696170458573605
$DArr[21,0] = 4.641614989579925
$DArr[22,0] = 4.444475428226083
$DArr[23,0] = 4.22297074739902
$ws.Range("D2:D25").Value = $DArr

$EArr = New-Object "object[,]" 24,1
$EArr[0,0] = 29.86504473444742
$EArr[1,0] = 27.69026846731998
$EArr[2,0] = 26.29205675092066
$EArr[3,0] = 25.70602877535457
$EArr[4,0] = 25.6076215683112
$EArr[5,0] = 26.28390721346732
$EArr[6,0] = 29.12759134308105
$EArr[7,0] = 34.22219387331137
$EArr[8,0] = 37.69621569722742
$EArr[9,0] = 39.22399010387477
$EArr[10,0] = 39.79597566759793
$EArr[11,0] = 39.67315116495583
$EArr[12,0] = 39.27121018502821
$EArr[13,0] = 39.02400292111651
$EArr[14,0] = 37.59429436235822
$EArr[15,0] = 36.70432070767297
$EArr[16,0] = 36.18782252526823
$EArr[17,0] = 36.01183757130035
$EArr[18,0] = 36.7995906127655
$EArr[19,0] = 39.38907693919166
$EArr[20,0] = 41.04221790569364
$EArr[21,0] = 40.16378053206288
$EArr[22,0] = 36.75721426407923
$EArr[23,0] = 32.89252398853505
$ws.Range("E2:E25").Value = $EArr

$FArr = New-Object "object[,]" 24,1
$FArr[0,0] = 17.75518195868449
$FArr[1,0] = 17.14176945398584
$FArr[2,0] = 16.77524196999751
$FArr[3,0] = 16.61400141833031
$FArr[4,0] = 16.56959065989545
$FArr[5,0] = 16.72497046929046
$FArr[6,0] = 17.48066969931257
$FArr[7,0] = 19.11861968476379
$FArr[8,0] = 20.34066971569571
$FArr[9,0] = 20.84195610334727
$FArr[10,0] = 21.06729681027591
$FArr[11,0] = 21.02545369730817
$FArr[12,0] = 20.86335585068473
$FArr[13,0] = 20.75075155625321
$FArr[14,0] = 20.19236493871767
$FArr[15,0] = 19.85176612296935
$FArr[16,0] = 19.69256271080111
$FArr[17,0] = 19.61110408100279
$FArr[18,0] = 19.89074133592759
$FArr[19,0] = 20.87301828555561
$FArr[20,0] = 21.56934619850209
$FArr[21,0] = 21.24243359141878
$FArr[22,0] = 19.94848042801264
$FArr[23,0] = 18.59025070022433
$ws.Range("F2:F25").Value = $FArr

$GArr = New-Object "object[,]" 24,1
$GArr[0,0] = 2.057930378311564
$GArr[1,0] = 2.064048891841666
$GArr[2,0] = 2.067904848094221
$GArr[3,0] = 2.069508935871258
$GArr[4,0] = 2.069785124282182
$GArr[5,0] = 2.067948371308099
$GArr[6,0] = 2.060047103224901
$GArr[7,0] = 2.045252384990889
$GArr[8,0] = 2.034760369430259
$GArr[9,0] = 2.030073816721957
$GArr[10,0] = 2.028291515289511
$GArr[11,0] = 2.028672275112224
$GArr[12,0] = 2.029925855351297
$GArr[13,0] = 2.030700154122882
$GArr[14,0] = 2.035116223417924
$GArr[15,0] = 2.037839544240324
$GArr[16,0] = 2.039396760793683
$GArr[17,0] = 2.039936579996586
$GArr[18,0] = 2.037547896057977
$GArr[19,0] = 2.029573276824686
$GArr[20,0] = 2.024380344840279
$GArr[21,0] = 2.027130046297772
$GArr[22,0] = 2.037648255556793
$GArr[23,0] = 2.049217081647693
$ws.Range("G2:G25").Value = $GArr

$IArr = New-Object "object[,]" 24,1
$IArr[0,0] = 3.607045034057078
$IArr[1,0] = 3.368673817656461
$IArr[2,0] = 3.219025104047763
$IArr[3,0] = 3.157456965053914
$IArr[4,0] = 3.147696317729529
$IArr[5,0] = 3.219281481872532
$IArr[6,0] = 3.526447503311346
$IArr[7,0] = 4.10380564727852
$IArr[8,0] = 4.516657697344205
$IArr[9,0] = 4.704272022066776
$IArr[10,0] = 4.774263864951045
$IArr[11,0] = 4.758984830741107
$IArr[12,0] = 4.709938037566803
$IArr[13,0] = 4.68036410262308
$IArr[14,0] = 4.507698184996459
$IArr[15,0] = 4.401505555253156
$IArr[16,0] = 4.338614581514044
$IArr[17,0] = 4.318543548889749
$IArr[18,0] = 4.412685807877819
$IArr[19,0] = 4.725745030086078
$IArr[20,0] = 4.927843456026763
$IArr[21,0] = 4.818625660600418
$IArr[22,0] = 4.404774556374785
$IArr[23,0] = 3.951725774018746
$ws.Range("I2:I25").Value = $IArr

$OArr = New-Object "object[,]" 24,1
$OArr[0,0] = 0
$OArr[1,0] = 0
$OArr[2,0] = 0
$OArr[3,0] = 0
$OArr[4,0] = 0
$OArr[5,0] = 0
$OArr[6,0] = 0
$OArr[7,0] = 0
$OArr[8,0] = 0
$OArr[9,0] = 0
$OArr[10,0] = 0
$OArr[11,0] = 0
$OArr[12,0] = 0
$OArr[13,0] = 0
$OArr[14,0] = 0
$OArr[15,0] = 0
$OArr[16,0] = 0
$OArr[17,0] = 0
$OArr[18,0] = 0
$OArr[19,0] = 0
$OArr[20,0] = 0
$OArr[21,0] = 0
$OArr[22,0] = 0
$OArr[23,0] = 0
$ws.Range("O2:O25").Value = $OArr

$PArr = New-Object "object[,]" 24,1
$PArr[0,0] = 0
$PArr[1,0] = 0
$PArr[2,0] = 0
$PArr[3,0] = 0
$PArr[4,0] = 0
$PArr[5,0] = 0
$PArr[6,0] = 0
$PArr[7,0] = 0
$PArr[8,0] = 0
$PArr[9,0] = 0
$PArr[10,0] = 0
$PArr[11,0] = 0
$PArr[12,0] = 0
$PArr[13,0] = 0
$PArr[14,0] = 0
$PArr[15,0] = 0
$PArr[16,0] = 0
$PArr[17,0] = 0
$PArr[18,0] = 0
$PArr[19,0] = 0
$PArr[20,0] = 0
$PArr[21,0] = 0
$PArr[22,0] = 0
$PArr[23,0] = 0
$ws.Range("P2:P25").Value = $PArr

$QArr = New-Object "object[,]" 24,1
$QArr[0,0] = 13.55683091514271
$QArr[1,0] = 13.33996087838208
$QArr[2,0] = 13.22219463795195
$QArr[3,0] = 13.16640053864074
$QArr[4,0] = 13.14325759494701
$QArr[5,0] = 13.18309357079213
$QArr[6,0] = 13.42905389519076
$QArr[7,0] = 14.10633865231042
$QArr[8,0] = 14.6715795628238
$QArr[9,0] = 14.89406522392085
$QArr[10,0] = 15.01102215433058
$QArr[11,0] = 14.99156269212319
$QArr[12,0] = 14.90613140823682
$QArr[13,0] = 14.84252412520621
$QArr[14,0] = 14.55769537661754
$QArr[15,0] = 14.38828333860699
$QArr[16,0] = 14.32352736784986
$QArr[17,0] = 14.27824197408994
$QArr[18,0] = 14.40851390499954
$QArr[19,0] = 14.89816928808926
$QArr[20,0] = 15.27773361539874
$QArr[21,0] = 15.11298383404167
$QArr[22,0] = 14.46401615914885
$QArr[23,0] = 13.84317342993667
$ws.Range("Q2:Q25").Value = $QArr
